$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) labels/order.
# New order: Year, totalAid, Pell Recipients, percentOfAidReceivers, AverageAid,
#            Percent of students with institutional aid, AverageInstitutionalAid
$ws.Range("A1").Value = "Year"
$ws.Range("B1").Value = "totalAid"
$ws.Range("C1").Value = "Pell Recipients"
$ws.Range("D1").Value = "percentOfAidReceivers"
$ws.Range("E1").Value = "AverageAid"
$ws.Range("F1").Value = "Percent of students with institutional aid"
$ws.Range("G1").Value = "AverageInstitutionalAid"

# Update the selected cell to G1 as in the target sheetView.
$ws.Range("G1").Select()
